$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 45.98144433333334
$ws.Range("H2").Value2 = 137.944333
$ws.Range("I2").Value2 = 0.9841234286873372
$ws.Range("J2").Value2 = 0.984123428687337
$ws.Range("M2").Value2 = 1.174933333333333
$ws.Range("N2").Value2 = 3.5248
$ws.Range("O2").Value2 = 0.01171850713626266
$ws.Range("P2").Value2 = 0.01171850713626266
$ws.Range("Q2").Value2 = 54.02513166204445
$ws.Range("R2").Value2 = 486.2261849584
$ws.Range("S2").Value2 = 0.01153245742203583
$ws.Range("T2").Value2 = 0.01153245742203583
$ws.Range("G3").Value2 = 45.98144433333334
$ws.Range("H3").Value2 = 137.944333
$ws.Range("I3").Value2 = 0.9841234286873372
$ws.Range("J3").Value2 = 0.984123428687337
$ws.Range("O3").Value2 = 0.2743256641287217
$ws.Range("P3").Value2 = 0.2743256641287218
$ws.Range("Q3").Value2 = 1264.707180744067
$ws.Range("R3").Value2 = 11382.3646266966
$ws.Range("S3").Value2 = 0.2699703131592885
$ws.Range("T3").Value2 = 0.2699703131592885
$ws.Range("G4").Value2 = 45.98144433333334
$ws.Range("H4").Value2 = 137.944333
$ws.Range("I4").Value2 = 0.9841234286873372
$ws.Range("J4").Value2 = 0.984123428687337
$ws.Range("M4").Value2 = 39.361408
$ws.Range("N4").Value2 = 118.084224
$ws.Range("O4").Value2 = 0.3925813724534833
$ws.Range("P4").Value2 = 0.3925813724534833
$ws.Range("Q4").Value2 = 1809.894390833622
$ws.Range("R4").Value2 = 16289.04951750259
$ws.Range("S4").Value2 = 0.3863485262977026
$ws.Range("T4").Value2 = 0.3863485262977025
$ws.Range("G5").Value2 = 45.98144433333334
$ws.Range("H5").Value2 = 137.944333
$ws.Range("I5").Value2 = 0.9841234286873372
$ws.Range("J5").Value2 = 0.984123428687337
$ws.Range("M5").Value2 = 32.221985
$ws.Range("N5").Value2 = 96.665955
$ws.Range("O5").Value2 = 0.3213744562815322
$ws.Range("P5").Value2 = 0.3213744562815322
$ws.Range("Q5").Value2 = 1481.613409587002
$ws.Range("R5").Value2 = 13334.52068628301
$ws.Range("S5").Value2 = 0.3162721318083102
$ws.Range("T5").Value2 = 0.3162721318083102
$ws.Range("I6").Value2 = 0.002244435796517234
$ws.Range("J6").Value2 = 0.002244435796517234
$ws.Range("M6").Value2 = 1.174933333333333
$ws.Range("N6").Value2 = 3.5248
$ws.Range("O6").Value2 = 0.01171850713626266
$ws.Range("P6").Value2 = 0.01171850713626266
$ws.Range("Q6").Value2 = 0.1232121255111111
$ws.Range("R6").Value2 = 1.1089091296
$ws.Range("S6").Value2 = 0.00002630143689837057
$ws.Range("T6").Value2 = 0.00002630143689837056
$ws.Range("I7").Value2 = 0.002244435796517234
$ws.Range("J7").Value2 = 0.002244435796517234
$ws.Range("O7").Value2 = 0.2743256641287217
$ws.Range("P7").Value2 = 0.2743256641287218
$ws.Range("S7").Value2 = 0.0006157063404738669
$ws.Range("T7").Value2 = 0.0006157063404738669
$ws.Range("I8").Value2 = 0.002244435796517234
$ws.Range("J8").Value2 = 0.002244435796517234
$ws.Range("M8").Value2 = 39.361408
$ws.Range("N8").Value2 = 118.084224
$ws.Range("O8").Value2 = 0.3925813724534833
$ws.Range("P8").Value2 = 0.3925813724534833
$ws.Range("Q8").Value2 = 4.127725893205334
$ws.Range("R8").Value2 = 37.149533038848
$ws.Range("S8").Value2 = 0.000881123685380463
$ws.Range("T8").Value2 = 0.0008811236853804628
$ws.Range("I9").Value2 = 0.002244435796517234
$ws.Range("J9").Value2 = 0.002244435796517234
$ws.Range("M9").Value2 = 32.221985
$ws.Range("N9").Value2 = 96.665955
$ws.Range("O9").Value2 = 0.3213744562815322
$ws.Range("P9").Value2 = 0.3213744562815322
$ws.Range("Q9").Value2 = 3.379033641656666
$ws.Range("R9").Value2 = 30.41130277491
$ws.Range("S9").Value2 = 0.0007213043337645338
$ws.Range("T9").Value2 = 0.0007213043337645338
$ws.Range("G10").Value2 = 0.547937
$ws.Range("H10").Value2 = 1.643811
$ws.Range("I10").Value2 = 0.01172728797372169
$ws.Range("J10").Value2 = 0.01172728797372169
$ws.Range("M10").Value2 = 1.174933333333333
$ws.Range("N10").Value2 = 3.5248
$ws.Range("O10").Value2 = 0.01171850713626266
$ws.Range("P10").Value2 = 0.01171850713626266
$ws.Range("Q10").Value2 = 0.6437894458666668
$ws.Range("R10").Value2 = 5.794105012799999
$ws.Range("S10").Value2 = 0.0001374263078090649
$ws.Range("T10").Value2 = 0.0001374263078090648
$ws.Range("G11").Value2 = 0.547937
$ws.Range("H11").Value2 = 1.643811
$ws.Range("I11").Value2 = 0.01172728797372169
$ws.Range("J11").Value2 = 0.01172728797372169
$ws.Range("O11").Value2 = 0.2743256641287217
$ws.Range("P11").Value2 = 0.2743256641287218
$ws.Range("Q11").Value2 = 15.07085887671867
$ws.Range("R11").Value2 = 135.637729890468
$ws.Range("S11").Value2 = 0.003217096061819975
$ws.Range("T11").Value2 = 0.003217096061819975
$ws.Range("G12").Value2 = 0.547937
$ws.Range("H12").Value2 = 1.643811
$ws.Range("I12").Value2 = 0.01172728797372169
$ws.Range("J12").Value2 = 0.01172728797372169
$ws.Range("M12").Value2 = 39.361408
$ws.Range("N12").Value2 = 118.084224
$ws.Range("O12").Value2 = 0.3925813724534833
$ws.Range("P12").Value2 = 0.3925813724534833
$ws.Range("Q12").Value2 = 21.567571815296
$ws.Range("R12").Value2 = 194.108146337664
$ws.Range("S12").Value2 = 0.004603914807880892
$ws.Range("T12").Value2 = 0.00460391480788089
$ws.Range("G13").Value2 = 0.547937
$ws.Range("H13").Value2 = 1.643811
$ws.Range("I13").Value2 = 0.01172728797372169
$ws.Range("J13").Value2 = 0.01172728797372169
$ws.Range("M13").Value2 = 32.221985
$ws.Range("N13").Value2 = 96.665955
$ws.Range("O13").Value2 = 0.3213744562815322
$ws.Range("P13").Value2 = 0.3213744562815322
$ws.Range("Q13").Value2 = 17.655617794945
$ws.Range("R13").Value2 = 158.900560154505
$ws.Range("S13").Value2 = 0.00376885079621176
$ws.Range("T13").Value2 = 0.00376885079621176
$ws.Range("E14").Value2 = 2
$ws.Range("F14").Value2 = 0.6666666666666666
$ws.Range("G14").Value2 = 0.08900066666666666
$ws.Range("H14").Value2 = 0.267002
$ws.Range("I14").Value2 = 0.001904847542424061
$ws.Range("J14").Value2 = 0.001904847542424061
$ws.Range("M14").Value2 = 1.174933333333333
$ws.Range("N14").Value2 = 3.5248
$ws.Range("O14").Value2 = 0.01171850713626266
$ws.Range("P14").Value2 = 0.01171850713626266
$ws.Range("Q14").Value2 = 0.1045698499555555
$ws.Range("R14").Value2 = 0.9411286495999999
$ws.Range("S14").Value2 = 0.00002232196951938875
$ws.Range("T14").Value2 = 0.00002232196951938874
$ws.Range("E15").Value2 = 2
$ws.Range("F15").Value2 = 0.6666666666666666
$ws.Range("G15").Value2 = 0.08900066666666666
$ws.Range("H15").Value2 = 0.267002
$ws.Range("I15").Value2 = 0.001904847542424061
$ws.Range("J15").Value2 = 0.001904847542424061
$ws.Range("O15").Value2 = 0.2743256641287217
$ws.Range("P15").Value2 = 0.2743256641287218
$ws.Range("Q15").Value2 = 2.447939247152889
$ws.Range("R15").Value2 = 22.031453224376
$ws.Range("S15").Value2 = 0.0005225485671394441
$ws.Range("T15").Value2 = 0.0005225485671394442
$ws.Range("E16").Value2 = 2
$ws.Range("F16").Value2 = 0.6666666666666666
$ws.Range("G16").Value2 = 0.08900066666666666
$ws.Range("H16").Value2 = 0.267002
$ws.Range("I16").Value2 = 0.001904847542424061
$ws.Range("J16").Value2 = 0.001904847542424061
$ws.Range("M16").Value2 = 39.361408
$ws.Range("N16").Value2 = 118.084224
$ws.Range("O16").Value2 = 0.3925813724534833
$ws.Range("P16").Value2 = 0.3925813724534833
$ws.Range("Q16").Value2 = 3.503191552938667
$ws.Range("R16").Value2 = 31.528723976448
$ws.Range("S16").Value2 = 0.0007478076625194829
$ws.Range("T16").Value2 = 0.0007478076625194828
$ws.Range("E17").Value2 = 2
$ws.Range("F17").Value2 = 0.6666666666666666
$ws.Range("G17").Value2 = 0.08900066666666666
$ws.Range("H17").Value2 = 0.267002
$ws.Range("I17").Value2 = 0.001904847542424061
$ws.Range("J17").Value2 = 0.001904847542424061
$ws.Range("M17").Value2 = 32.221985
$ws.Range("N17").Value2 = 96.665955
$ws.Range("O17").Value2 = 0.3213744562815322
$ws.Range("P17").Value2 = 0.3213744562815322
$ws.Range("Q17").Value2 = 2.867778146323333
$ws.Range("R17").Value2 = 25.81000331691
$ws.Range("S17").Value2 = 0.0006121693432457455
$ws.Range("T17").Value2 = 0.0006121693432457456
